$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row labels: exponent changes from 10^{-3} to 10^{3}
# (the visible text / cell value is the only real change; Excel will
# reorganize the shared-strings table on save)
$ws.Range("A1").Value = " `$x / 10^{3} \unit\meter`$"
$ws.Range("B1").Value = "`$D_0 (x) / 10^{3} \unit\meter`$"
$ws.Range("C1").Value = "`$D_m (x)/ 10^{3} \unit\meter`$"

# Move the active selection to F8, matching the author's editing position
$ws.Range("F8").Select()
